$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 72 (this shifts the old rows 72..157 down to 73..158)
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new weekly record
$ws.Range("A72").Value = 8
$ws.Range("B72").Value = "Terminal La Palmera de La Serena"
$ws.Range("C72").Value = "Coquimbo"
$ws.Range("D72").Value = 44763
$ws.Range("E72").Value = 4
$ws.Range("F72").Value = 100112040
$ws.Range("G72").Value = "Cilantro"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 2000
$ws.Range("K72").Value = 2000
$ws.Range("L72").Value = 2500
$ws.Range("M72").Value = 2250
$ws.Range("N72").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O72").Value = "Provincia del Elquí"
$ws.Range("P72").Value = 1500
$ws.Range("Q72").Value = 1.5
$ws.Range("R72").Value = "Hortaliza"
